$d = $word.ActiveDocument

# The Bibliografia paragraph currently holds one run with all references
# concatenated together. We split it into multiple <w:t> runs separated by
# manual line breaks (<w:br/>), matching each reference on its own line.
# We do this by finding unique short anchors right at each split point and
# replacing them with themselves plus a manual line break ("^l").

$splits = @(
    @{ Find = "2022."; Replace = "2022.^l" },
    @{ Find = "2019."; Replace = "2019.^l" },
    @{ Find = "Novatec, 2015Downey"; Replace = "Novatec, 2015^lDowney" },
    @{ Find = "2016."; Replace = "2016.^l" },
    @{ Find = "2014."; Replace = "2014.^l" },
    @{ Find = "2008."; Replace = "2008.^l" },
    @{ Find = "2006."; Replace = "2006.^l" }
)

foreach ($pair in $splits) {
    $rng = $d.Content
    $find = $rng.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2) | Out-Null
}
